# Update workbook "广州-漫展信息.xlsx" to the newly scraped data snapshot.
# - Refreshes "想去人数" (interest-count) values across all sheets.
# - Removes one duplicate "Look Look动漫嘉年华" row from 展览 (sheet1) and
#   全部类型 (sheet4), shifting the two rows below it up and dropping the
#   now-empty trailing row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (exhibitions)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value  = 8320
$ws1.Range("F3").Value  = 133
$ws1.Range("F4").Value  = 104
$ws1.Range("F5").Value  = 35801
$ws1.Range("F9").Value  = 466
$ws1.Range("F13").Value = 69
$ws1.Range("F14").Value = 638
$ws1.Range("F15").Value = 452
$ws1.Range("F17").Value = 586
$ws1.Range("F19").Value = 434
$ws1.Range("F20").Value = 430
$ws1.Range("F21").Value = 1128
$ws1.Range("F23").Value = 746
$ws1.Range("F24").Value = 2405
$ws1.Range("F25").Value = 888
$ws1.Range("F26").Value = 511
$ws1.Range("F28").Value = 1105
$ws1.Range("F30").Value = 680

# Row 31 becomes the old row 32's event (A/B stay as-is)
$ws1.Range("C31").Value = "广州·第五届AP动漫嘉年华"
$ws1.Range("D31").Value = "西环路1号 广州岭南会展中心"
$ws1.Range("E31").Value = "2024.06.01 10:00-06.01 17:00"
$ws1.Range("F31").Value = 16
$ws1.Range("G31").Value = 55
$ws1.Range("H31").Value = "https://show.bilibili.com/platform/detail.html?id=83462"
$ws1.Range("I31").Value = "//i1.hdslb.com/bfs/openplatform/202403/ZR2jKMOg1711076939687.jpeg"

# Row 32 becomes the old row 33's event (A stays as-is)
$ws1.Range("B32").Value = "2024-06-22"
$ws1.Range("C32").Value = "广州·622排球少年only"
$ws1.Range("D32").Value = "岭南购物城内 广州OMG网红街"
$ws1.Range("E32").Value = "2024.06.22 10:00-06.22 17:30"
$ws1.Range("F32").Value = 1104
$ws1.Range("G32").Value = 68
$ws1.Range("H32").Value = "https://show.bilibili.com/platform/detail.html?id=82974"
$ws1.Range("I32").Value = "//i0.hdslb.com/bfs/openplatform/202403/WMlOXSZn1710748067155.jpeg"

# The old row 33 no longer exists once the duplicate is gone.
$ws1.Rows.Item(33).Delete()

# ---------------------------------------------------------------------------
# Sheet "演出" (performances)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 357
$ws2.Range("F9").Value = 139

# ---------------------------------------------------------------------------
# Sheet "本地生活" (local life)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 564

# ---------------------------------------------------------------------------
# Sheet "全部类型" (all types combined) - mirrors the 展览 sheet edits above.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value  = 564
$ws4.Range("F3").Value  = 8320
$ws4.Range("F4").Value  = 133
$ws4.Range("F5").Value  = 104
$ws4.Range("F7").Value  = 35801
$ws4.Range("F11").Value = 466
$ws4.Range("F15").Value = 357
$ws4.Range("F19").Value = 69
$ws4.Range("F20").Value = 638
$ws4.Range("F21").Value = 452
$ws4.Range("F25").Value = 139
$ws4.Range("F28").Value = 586
$ws4.Range("F30").Value = 434
$ws4.Range("F31").Value = 430
$ws4.Range("F32").Value = 1128
$ws4.Range("F34").Value = 746
$ws4.Range("F35").Value = 2405
$ws4.Range("F36").Value = 888
$ws4.Range("F37").Value = 511
$ws4.Range("F39").Value = 1105
$ws4.Range("F42").Value = 680

# Row 43 becomes the old row 44's event (A/B stay as-is)
$ws4.Range("C43").Value = "广州·第五届AP动漫嘉年华"
$ws4.Range("D43").Value = "西环路1号 广州岭南会展中心"
$ws4.Range("E43").Value = "2024.06.01 10:00-06.01 17:00"
$ws4.Range("F43").Value = 16
$ws4.Range("G43").Value = 55
$ws4.Range("H43").Value = "https://show.bilibili.com/platform/detail.html?id=83462"
$ws4.Range("I43").Value = "//i1.hdslb.com/bfs/openplatform/202403/ZR2jKMOg1711076939687.jpeg"

# Row 44 becomes the old row 45's event (A stays as-is)
$ws4.Range("B44").Value = "2024-06-22"
$ws4.Range("C44").Value = "广州·622排球少年only"
$ws4.Range("D44").Value = "岭南购物城内 广州OMG网红街"
$ws4.Range("E44").Value = "2024.06.22 10:00-06.22 17:30"
$ws4.Range("F44").Value = 1104
$ws4.Range("G44").Value = 68
$ws4.Range("H44").Value = "https://show.bilibili.com/platform/detail.html?id=82974"
$ws4.Range("I44").Value = "//i0.hdslb.com/bfs/openplatform/202403/WMlOXSZn1710748067155.jpeg"

# The old row 45 no longer exists once the duplicate is gone.
$ws4.Rows.Item(45).Delete()
